$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure the phone-number column keeps its text formatting (so values such as
# "0987654321" keep their leading zero instead of being auto-converted to a
# number) and then reset the style back to Normal so no extra cell style
# index gets stamped onto the cells.
$phoneCells = "C6:C10"
$ws.Range($phoneCells).NumberFormat = "@"

# --- Row 5 (customer #4) ---
$ws.Range("B5").Value = "Phú"
$ws.Range("C5").Value = 1234567890
$ws.Range("D5").Value = "phu@edu.vn"
$ws.Range("E5").Value = "Tây Ninh"

# --- Row 6 (customer #5) ---
$ws.Range("B6").Value = "Phú nè"
$ws.Range("C6").Value = "0987654321"
$ws.Range("D6").Value = "plplpl@mail.vn"
$ws.Range("E6").Value = "Nhà"

# --- Row 7 (customer #6, was labelled 7, now 6) ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Phú"
$ws.Range("C7").Value = "0987654321"
$ws.Range("D7").Value = "mail@gmail.com"
$ws.Range("E7").Value = "Tây Ninh"

# --- Row 8 (customer #7, new row) ---
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Phú"
$ws.Range("C8").Value = "0987654321"
$ws.Range("D8").Value = "mail@mail.com"
$ws.Range("E8").Value = "Tây Ninh"

# --- Row 9 (customer #8, new row) ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "hehe"
$ws.Range("C9").Value = "0987654321"
$ws.Range("D9").Value = "m@m.com"
$ws.Range("E9").Value = "tn"

# --- Row 10 (customer #9, new row) ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Hiếu gà"
$ws.Range("C10").Value = "0987654321"
$ws.Range("D10").Value = "hieuga2003@gmail.com"
$ws.Range("E10").Value = "Lã Xuân Oai"

# Restore the default "Normal" style on the phone-number cells so they don't
# keep a custom style index (only the text cell type should differ).
$ws.Range($phoneCells).Style = "Normal"
